$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value for every data row (2..526).
# All of these values change from 45202 (2023-10-03) to 45203 (2023-10-04).
$ws.Range("C2:C526").Value = 45203
